# Buttons.pptx edit:
#  1. Title slide (slide 1): "Buttons & Functions" -> "Buttons"
#  2. Delete the "Functions recap" slide (slide 9). This removes that slide
#     (and its speaker notes) entirely; the following "Button example" slide
#     shifts up to become slide 9, and "Questions?" shifts up to become slide 10.

$p = $ppt.ActivePresentation

# 1. Update the title text on the first slide.
$titleShape = $p.Slides.Item(1).Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "Buttons"

# 2. Delete the "Functions recap" slide (currently slide 9).
$p.Slides.Item(9).Delete()
